$wb = $excel.ActiveWorkbook

# --- Goals sheet: rotate Goal IDs to new UUIDs ---
$goalsWs = $wb.Worksheets.Item("Goals")
$goalsWs.Range("A2").Value = "0964ea5a-ea5e-4c45-9787-1479b64fa2ca"
$goalsWs.Range("A3").Value = "aa0ec478-ef74-4b20-b014-e9487258c775"
$goalsWs.Range("A4").Value = "b21c8691-eb0f-4674-97fd-78b6d2be89ae"
$goalsWs.Range("A5").Value = "b32b6b14-36d0-492e-9353-e838b3ce8f99"

# --- Tasks sheet: update Task Status (and completion flag) ---
$tasksWs = $wb.Worksheets.Item("Tasks")

# Row 2: [Week1] Learning English Speaking -> "In Progress"
$tasksWs.Range("H2").Value = "In Progress"

# Row 7: Test 2 (Learn Java) -> "To Do"
$tasksWs.Range("H7").Value = "To Do"

# Row 12: Improve Code (Learn Spring Boot) -> "Done", now completed
$tasksWs.Range("H12").Value = "Done"
$tasksWs.Range("I12").Value = $true
